# Rewrite the match-log table: insert "ownTeam"/"oppTeam" columns after
# "result", reorder the rows chronologically, and refresh every value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@("venue", "date", "result", "ownTeam", "oppTeam", "batsman", "totalRuns", "totalBalls", "total4s", "total6s", "sr")
  ,@(" Abu Dhabi", " October 16 2020", "Mumbai won by 8 wickets (with 19 balls remaining)", "Kolkata Knight Riders", "Mumbai Indians", "Dinesh Karthik †", "4", "8", "1", "0", "50.00")
  ,@(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Kolkata Knight Riders", "Sunrisers Hyderabad", "Dinesh Karthik †", "29", "14", "2", "2", "207.14")
  ,@(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Dinesh Karthik †", "0", "1", "0", "0", "0.00")
  ,@(" Dubai (DSC)", " October 29 2020", "Super Kings won by 6 wickets", "Kolkata Knight Riders", "Chennai Super Kings", "Dinesh Karthik †", "21", "10", "3", "0", "210.00")
  ,@(" Sharjah", " October 26 2020", "Kings XI won by 8 wickets (with 7 balls remaining)", "Kolkata Knight Riders", "Kings XI Punjab", "Dinesh Karthik †", "0", "2", "0", "0", "0.00")
  ,@(" Abu Dhabi", " October 21 2020", "RCB won by 8 wickets (with 39 balls remaining)", "Kolkata Knight Riders", "Royal Challengers Bangalore", "Dinesh Karthik †", "4", "14", "0", "0", "28.57")
  ,@(" Abu Dhabi", " October 24 2020", "KKR won by 59 runs", "Kolkata Knight Riders", "Delhi Capitals", "Dinesh Karthik †", "3", "6", "0", "0", "50.00")
)

$rowCount = $data.Count
$colCount = $data[0].Count

# Clear the previously-used range first (old layout was only 9 columns wide).
$ws.Cells.Clear()

# Target range, formatted as Text first so numeric-looking strings (e.g.
# "4", "50.00") round-trip as text instead of being coerced to numbers.
$fullRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($rowCount, $colCount))
$fullRange.NumberFormat = "@"

for ($r = 0; $r -lt $rowCount; $r++) {
  for ($c = 0; $c -lt $colCount; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
  }
}
